$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New feed rows appended by the workflow run.
$rows = @(
    @{ Link = "https://www.genomeweb.com/cancer/amoydx-lung-cancer-panel-nabs-japanese-approval-cdx-non-small-cell-lung-cancer-drug-ibtrozi"; Keywords = "CDx"; Title = "AmoyDx Lung Cancer Panel Nabs Japanese Approval as CDx for Non-Small Cell Lung Cancer Drug Ibtrozi" },
    @{ Link = "https://www.360dx.com/cancer/amoydx-lung-cancer-panel-nabs-japanese-approval-cdx-non-small-cell-lung-cancer-drug-ibtrozi"; Keywords = "CDx"; Title = "AmoyDx Lung Cancer Panel Nabs Japanese Approval as CDx for Non-Small Cell Lung Cancer Drug Ibtrozi" },
    @{ Link = "https://www.medpagetoday.com/meetingcoverage/ims/117536"; Keywords = "BCMA"; Title = "Elranatamab-Based Triplet Effective in Transplant-Ineligible Myeloma" }
)

$startRow = 43
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    $ws.Cells.Item($r, 2).Value = $row.Keywords
    $ws.Cells.Item($r, 3).Value = $row.Title

    $cell = $ws.Cells.Item($r, 1)
    $ws.Hyperlinks.Add($cell, $row.Link) | Out-Null
    $cell.Style = $ws.Cells.Item($r - 1, 1).Style
}

"done"
